# Fix bullets/lists indentation in the numbering definitions:
# for every w:lvl in the abstractNum, drop the <w:tabs><w:tab w:val="num" .../></w:tabs>
# element from w:pPr and bump w:ind/@w:left by 240 twips (except level 0, which
# stays at 480 since its tab stop was already at 0).

$d = $word.ActiveDocument

# Tab stop position (at level N) -> (old w:ind left, new w:ind left)
$tabPositions = @(0, 720, 1440, 2160, 2880, 3600, 4320, 5040, 5760)
$oldLefts     = @(480, 1200, 1920, 2640, 3360, 4080, 4800, 5520, 6240)
$newLefts     = @(480, 1440, 2160, 2880, 3600, 4320, 5040, 5760, 6480)

$xml = $d.Content.WordOpenXML

for ($i = 0; $i -lt $tabPositions.Length; $i++) {
    $tabPos  = $tabPositions[$i]
    $oldLeft = $oldLefts[$i]
    $newLeft = $newLefts[$i]

    $oldBlock = '<w:tabs><w:tab w:val="num" w:pos="' + $tabPos + '" /></w:tabs><w:ind w:left="' + $oldLeft + '" w:hanging="480" />'
    $newBlock = '<w:ind w:left="' + $newLeft + '" w:hanging="480" />'

    if ($xml.IndexOf($oldBlock) -lt 0) {
        Write-Output "WARNING: pattern not found for tab pos $tabPos"
    }

    $xml = $xml.Replace($oldBlock, $newBlock)
}

$d.Content.WordOpenXML = $xml

Write-Output "done"
